$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) First three rows: change existing values to "0M"
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# 2) Insert 10 new single-value rows right after row 3 (before the current row 4)
# Rows.Add always inserts directly before the same reference row, so pushing
# values in from last-to-first yields the correct final top-to-bottom order.
$newValues = @("99", "0.00002", "0.00005", "0.00003", "0.00001", "0.00003", "0.00003", "0.00004", "0.00327", "100.0")

$beforeRow = $t.Rows.Item(4)
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($beforeRow)
    $t.Cell($newRow.Index, 1).Range.Text = $newValues[$i]
}

# 3) Collapse the final three multi-run/tab rows down to single simple values
$rowCount = $t.Rows.Count
$t.Cell($rowCount - 2, 1).Range.Text = "99.99"
$t.Cell($rowCount - 1, 1).Range.Text = "0"
$t.Cell($rowCount, 1).Range.Text = "27"
